$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the promo-code / Amazon-link pairs for rows 4-7 first, then 2-3,
# so new shared-string entries are appended in the same relative order as
# the target workbook's sharedStrings table.
$ws.Range("A4").Value = "75E5GYPF"
$ws.Range("B4").Value = "https://amzn.to/35YdcNT"

$ws.Range("A5").Value = "50G2LSZS"
$ws.Range("B5").Value = "https://amzn.to/34SCPk3"

$ws.Range("A6").Value = "70P4RNKQ"
$ws.Range("B6").Value = "https://amzn.to/2Gu3zOo"

$ws.Range("A7").Value = "5NJ7BSJG"
$ws.Range("B7").Value = "https://amzn.to/3kQbTa8"

$ws.Range("A2").Value = "69226XKI"
$ws.Range("B2").Value = "https://amzn.to/2GphbKL"

$ws.Range("A3").Value = "5049YPD2"
$ws.Range("B3").Value = "https://amzn.to/3enAIYD"
